$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1853106.9
$ws.Range("J17").Value = 1853106.9
$ws.Range("L17").Value = 5559320.699999999
$ws.Range("N17").Value = -5559656.699999999
$ws.Range("H33").Value = 155.8077
$ws.Range("I33").Value = 114.05556
$ws.Range("K33").Value = 114.05556
$ws.Range("M33").Value = 114.94444
$ws.Range("H64").Value = 3808.225
$ws.Range("I64").Value = 3591.258
$ws.Range("J64").Value = 4555.5557
$ws.Range("K64").Value = 3591.258
$ws.Range("L64").Value = 4555.5557
$ws.Range("M64").Value = -3343.258
$ws.Range("N64").Value = -5051.5557
$ws.Range("H67").Value = 3808.225
$ws.Range("I67").Value = 3591.258
$ws.Range("J67").Value = 4555.5557
$ws.Range("K67").Value = 3591.258
$ws.Range("L67").Value = 4555.5557
$ws.Range("M67").Value = -2733.258
$ws.Range("N67").Value = -6271.5557
$ws.Range("H76").Value = 3190.6592
$ws.Range("I76").Value = 2571.9656
$ws.Range("K76").Value = 2571.9656
$ws.Range("M76").Value = -2256.9656
$ws.Range("H79").Value = 3190.6592
$ws.Range("I79").Value = 2571.9656
$ws.Range("K79").Value = 2571.9656
$ws.Range("M79").Value = -1479.9656
$ws.Range("H87").Value = 13881.287
$ws.Range("J87").Value = 13881.287
$ws.Range("L87").Value = 13881.287
$ws.Range("N87").Value = -16377.287
$ws.Range("H90").Value = 13881.287
$ws.Range("J90").Value = 13881.287
$ws.Range("L90").Value = 41643.861
$ws.Range("N90").Value = -54123.861
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1262.1333
$ws.Range("I2").Value = 1181.8334
$ws.Range("J2").Value = 1583.3334
$ws.Range("K2").Value = 1181.8334
$ws.Range("L2").Value = 1583.3334
$ws.Range("M2").Value = -1068.8334
$ws.Range("N2").Value = -1809.3334
$ws.Range("H63").Value = 3625.6667
$ws.Range("I63").Value = 3598.9285
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3598.9285
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -2912.9285
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 3625.6667
$ws.Range("I66").Value = 3598.9285
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 17994.6425
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -14562.6425
$ws.Range("N66").Value = -26864
$ws.Range("H88").Value = 1798.091
$ws.Range("I88").Value = 1396.6666
$ws.Range("J88").Value = 2279.8
$ws.Range("K88").Value = 1396.6666
$ws.Range("L88").Value = 2279.8
$ws.Range("M88").Value = -990.6666
$ws.Range("N88").Value = -3091.8
$ws.Range("H91").Value = 1798.091
$ws.Range("I91").Value = 1396.6666
$ws.Range("J91").Value = 2279.8
$ws.Range("K91").Value = 1396.6666
$ws.Range("L91").Value = 2279.8
$ws.Range("M91").Value = 7.333399999999983
$ws.Range("N91").Value = -5087.8
$ws.Range("H116").Value = 1262.1333
$ws.Range("I116").Value = 1181.8334
$ws.Range("J116").Value = 1583.3334
$ws.Range("K116").Value = 1181.8334
$ws.Range("L116").Value = 1583.3334
$ws.Range("M116").Value = 1112.1666
$ws.Range("N116").Value = -6171.3334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1262.1333
$ws.Range("I3").Value = 1181.8334
$ws.Range("J3").Value = 1583.3334
$ws.Range("K3").Value = 1181.8334
$ws.Range("L3").Value = 1583.3334
$ws.Range("M3").Value = -1067.8334
$ws.Range("N3").Value = -1811.3334
$ws.Range("H86").Value = 2114.6843
$ws.Range("I86").Value = 2031.1875
$ws.Range("J86").Value = 2560
$ws.Range("K86").Value = 2031.1875
$ws.Range("L86").Value = 2560
$ws.Range("M86").Value = -908.1875
$ws.Range("N86").Value = -4806
$ws.Range("H89").Value = 2114.6843
$ws.Range("I89").Value = 2031.1875
$ws.Range("J89").Value = 2560
$ws.Range("K89").Value = 10155.9375
$ws.Range("L89").Value = 12800
$ws.Range("M89").Value = -4539.9375
$ws.Range("N89").Value = -24032
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5983.727
$ws.Range("I62").Value = 3979.4443
$ws.Range("J62").Value = 15003
$ws.Range("K62").Value = 3979.4443
$ws.Range("L62").Value = 15003
$ws.Range("M62").Value = -3355.4443
$ws.Range("N62").Value = -16251
$ws.Range("H65").Value = 5983.727
$ws.Range("I65").Value = 3979.4443
$ws.Range("J65").Value = 15003
$ws.Range("K65").Value = 19897.2215
$ws.Range("L65").Value = 75015
$ws.Range("M65").Value = -16777.2215
$ws.Range("N65").Value = -81255
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H135").Value = 42097.5
$ws.Range("J135").Value = 42097.5
$ws.Range("L135").Value = 42097.5
$ws.Range("N135").Value = -52237.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5127.827
$ws.Range("I70").Value = 4431.8276
$ws.Range("J70").Value = 6005.391
$ws.Range("K70").Value = 4431.8276
$ws.Range("L70").Value = 6005.391
$ws.Range("M70").Value = -4161.8276
$ws.Range("N70").Value = -6545.391
$ws.Range("H73").Value = 5127.827
$ws.Range("I73").Value = 4431.8276
$ws.Range("J73").Value = 6005.391
$ws.Range("K73").Value = 4431.8276
$ws.Range("L73").Value = 6005.391
$ws.Range("M73").Value = -3495.8276
$ws.Range("N73").Value = -7877.391
$ws.Range("H80").Value = 2761.5386
$ws.Range("I80").Value = 2640
$ws.Range("J80").Value = 3166.6667
$ws.Range("K80").Value = 2640
$ws.Range("L80").Value = 3166.6667
$ws.Range("M80").Value = -1642
$ws.Range("N80").Value = -5162.6667
$ws.Range("H83").Value = 2761.5386
$ws.Range("I83").Value = 2640
$ws.Range("J83").Value = 3166.6667
$ws.Range("K83").Value = 13200
$ws.Range("L83").Value = 15833.3335
$ws.Range("M83").Value = -8208
$ws.Range("N83").Value = -25817.3335
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1757.5883
$ws.Range("I82").Value = 1104.6
$ws.Range("K82").Value = 1104.6
$ws.Range("M82").Value = -743.5999999999999
$ws.Range("H85").Value = 1757.5883
$ws.Range("I85").Value = 1104.6
$ws.Range("K85").Value = 1104.6
$ws.Range("M85").Value = 143.4000000000001
$ws.Range("H122").Value = 12348166
$ws.Range("I122").Value = 22223384
$ws.Range("K122").Value = 66670152
$ws.Range("M122").Value = -66667702
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5851174.5
$ws.Range("I136").Value = 10417218
$ws.Range("J136").Value = 6639.2
$ws.Range("K136").Value = 31251654
$ws.Range("L136").Value = 19917.6
$ws.Range("M136").Value = -31249104
$ws.Range("N136").Value = -25017.6
